# feat: add 2022-Q4 data
#
# 1. Insert a new worksheet named "2022-Q4" right after "总计" (so the
#    sheet order becomes: 总计, 2022-Q4, 2022-Q3, 2022-Q1, 2021-Q2) and
#    fill it with the new fund-holding data for that quarter.
# 2. Update the "总计" (summary) sheet: insert a new row for 2022-Q4 right
#    after the header, which pushes the existing 2022-Q3 / 2022-Q1 /
#    2021-Q2 rows down by one, and fix up the "持有数量(只)" /
#    "持有市值(亿元)" numbers to match the new totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: add the new "2022-Q4" sheet right after "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Add($null, $totalSheet)
$q4.Name = "2022-Q4"

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q4.Cells.Item(1, $c + 2).Value = $headers[$c]
}

$q4data = @(
    @("014133", "工银中证500六个月持有指数增强A", "1.59", "94.29", "1.28", "0.0204", 7),
    @("014134", "工银中证500六个月持有指数增强C", "0.88", "94.29", "1.28", "0.0113", 7)
)

for ($r = 0; $r -lt $q4data.Length; $r++) {
    $row = $q4data[$r]

    $q4.Cells.Item($r + 2, 1).Value = $r

    for ($c = 0; $c -lt $row.Length; $c++) {
        $cell = $q4.Cells.Item($r + 2, $c + 2)
        if ($c -lt 6) {
            # keep fund code / name / ratios as text, like the other sheets
            $cell.NumberFormat = "@"
        }
        $cell.Value = $row[$c]
    }
}

# ---------------------------------------------------------------------
# Step 2: update the "总计" summary sheet with the new 2022-Q4 row and
# shift the following rows down.
# ---------------------------------------------------------------------
$totalData = @(
    @(0, "2022-Q4", 2, 0.03),
    @(1, "2022-Q3", 2, 0.02),
    @(2, "2022-Q1", 3, 0.13),
    @(3, "2021-Q2", 2, 0.01)
)

for ($r = 0; $r -lt $totalData.Length; $r++) {
    $row = $totalData[$r]
    $totalSheet.Cells.Item($r + 2, 1).Value = $row[0]
    $totalSheet.Cells.Item($r + 2, 2).Value = $row[1]
    $totalSheet.Cells.Item($r + 2, 3).Value = $row[2]
    $totalSheet.Cells.Item($r + 2, 4).Value = $row[3]
}

Write-Host "edit complete"
